$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$ws.Range("F2").Value = 3.677627633308154
$ws.Range("I2").Value = 57.96693520079926
$ws.Range("C3").Value = 3.151019696039055
$ws.Range("F3").Value = 7.238308117577713
$ws.Range("C4").Value = 0.2646268688594923
$ws.Range("F4").Value = 5.740000133751053
$ws.Range("C5").Value = 0.2619810684096068
$ws.Range("F5").Value = 5.498352708913107
$ws.Range("F6").Value = 7.160682698745746
$ws.Range("F7").Value = 7.160431964469049
$ws.Range("F8").Value = 7.160410566463508
$ws.Range("F9").Value = 7.170478957423848
$ws.Range("C10").Value = 2.007575984369963
$ws.Range("F10").Value = 7.16064242014708
$ws.Range("C11").Value = 2.079712437168695
$ws.Range("C15").Value = 3.150995277138613
$ws.Range("C16").Value = 2.272233565660938
$ws.Range("C17").Value = 0.3167712909515016
$ws.Range("C21").Value = 3.151041094044596
$ws.Range("C23").Value = 0.3474650935825892
$ws.Range("C27").Value = 3.150995277138613
$ws.Range("C28").Value = 1.682770875173248
$ws.Range("C29").Value = 2.326916546433885
$ws.Range("C33").Value = 3.151041094044596
$ws.Range("C34").Value = 1.682474323990569
$ws.Range("C35").Value = 2.326916546433885
$ws.Range("C39").Value = 3.151019696039055
$ws.Range("C40").Value = 1.682474323990569
$ws.Range("C41").Value = 2.326916546433885
$ws.Range("C46").Value = 1.682770875173248
$ws.Range("C51").Value = 3.151019696039055
$ws.Range("C52").Value = 1.682474323990569
$ws.Range("C53").Value = 2.327148400117457
